$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4880.769
$ws.Range("J40").Value = 4825
$ws.Range("L40").Value = 4825
$ws.Range("N40").Value = -5175

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 666.3570999999999
$ws.Range("J70").Value = 983.1
$ws.Range("L70").Value = 2949.3
$ws.Range("N70").Value = -3489.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 666.3570999999999
$ws.Range("J73").Value = 983.1
$ws.Range("L73").Value = 2949.3
$ws.Range("N73").Value = -4821.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 8000
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 8000
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1412
$ws.Range("I94").Value = 1412
$ws.Range("K94").Value = 1412
$ws.Range("M94").Value = -961

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2747.913
$ws.Range("J138").Value = 4233.7144
$ws.Range("L138").Value = 12701.1432
$ws.Range("N138").Value = -22981.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13243
$ws.Range("I2").Value = 17060.4
$ws.Range("K2").Value = 17060.4
$ws.Range("M2").Value = -16947.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3132.35
$ws.Range("I45").Value = 2874.8572
$ws.Range("K45").Value = 2874.8572
$ws.Range("M45").Value = -2497.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6288.7856
$ws.Range("I61").Value = 6288.7856
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6288.7856
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -6076.7856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1791.0834
$ws.Range("I74").Value = 1946.5
$ws.Range("K74").Value = 1946.5
$ws.Range("M74").Value = -1072.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1791.0834
$ws.Range("I77").Value = 1946.5
$ws.Range("K77").Value = 9732.5
$ws.Range("M77").Value = -5364.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 200602
$ws.Range("J101").Value = 200602
$ws.Range("L101").Value = 200602
$ws.Range("N101").Value = -207092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3436.2
$ws.Range("I102").Value = 3546.8572
$ws.Range("J102").Value = 1887
$ws.Range("K102").Value = 3546.8572
$ws.Range("L102").Value = 1887
$ws.Range("M102").Value = -1924.8572
$ws.Range("N102").Value = -5131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 13243
$ws.Range("I116").Value = 17060.4
$ws.Range("K116").Value = 17060.4
$ws.Range("M116").Value = -14766.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3416.745
$ws.Range("I122").Value = 3141.818
$ws.Range("J122").Value = 5144.857
$ws.Range("K122").Value = 9425.454000000002
$ws.Range("L122").Value = 15434.571
$ws.Range("M122").Value = -6975.454000000002
$ws.Range("N122").Value = -20334.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3444.5454
$ws.Range("I132").Value = 3439.15
$ws.Range("K132").Value = 10317.45
$ws.Range("M132").Value = -7787.450000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6288.7856
$ws.Range("I136").Value = 6288.7856
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 18866.3568
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -16316.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13243
$ws.Range("I3").Value = 17060.4
$ws.Range("K3").Value = 17060.4
$ws.Range("M3").Value = -16946.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 35753.832
$ws.Range("I99").Value = 58078
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 58078
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -56580
$ws.Range("N99").Value = -7496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2476.5789
$ws.Range("I105").Value = 858.1539
$ws.Range("K105").Value = 858.1539
$ws.Range("M105").Value = 888.8461

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7988.7393
$ws.Range("I107").Value = 7564.684
$ws.Range("K107").Value = 7564.684
$ws.Range("M107").Value = -5644.684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1340.2354
$ws.Range("I31").Value = 1285.4584
$ws.Range("J31").Value = 2216.6667
$ws.Range("K31").Value = 1285.4584
$ws.Range("L31").Value = 2216.6667
$ws.Range("M31").Value = -990.4584
$ws.Range("N31").Value = -2806.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1340.2354
$ws.Range("I34").Value = 1285.4584
$ws.Range("J34").Value = 2216.6667
$ws.Range("K34").Value = 1285.4584
$ws.Range("L34").Value = 2216.6667
$ws.Range("M34").Value = -1083.4584
$ws.Range("N34").Value = -2620.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3898.3235
$ws.Range("I107").Value = 801.86957
$ws.Range("K107").Value = 801.86957
$ws.Range("M107").Value = 1118.13043

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 10096.5
$ws.Range("J12").Value = 11725.8
$ws.Range("L12").Value = 35177.39999999999
$ws.Range("N12").Value = -35523.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 41667270
$ws.Range("I29").Value = 1046.5
$ws.Range("J29").Value = 83333500
$ws.Range("K29").Value = 3139.5
$ws.Range("L29").Value = 250000500
$ws.Range("M29").Value = -2862.5
$ws.Range("N29").Value = -250001054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 459.5
$ws.Range("I31").Value = 419
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 1257
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -969
$ws.Range("N31").Value = -2076

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4133.3335
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 12400.0005
$ws.Range("N94").Value = -13752.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2494.5334
$ws.Range("I137").Value = 2371.7273
$ws.Range("K137").Value = 7115.1819
$ws.Range("M137").Value = -2015.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4219.421
$ws.Range("J102").Value = 5190
$ws.Range("L102").Value = 5190
$ws.Range("N102").Value = -8434

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6473.8335
$ws.Range("I113").Value = 4703.8
$ws.Range("K113").Value = 4703.8
$ws.Range("M113").Value = -2533.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2668.6956
$ws.Range("I122").Value = 2049.0557
$ws.Range("K122").Value = 6147.1671
$ws.Range("M122").Value = -3697.1671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 463.2857
$ws.Range("I55").Value = 470.29166
$ws.Range("K55").Value = 470.29166
$ws.Range("M55").Value = -297.29166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3039
$ws.Range("I122").Value = 3118.875
$ws.Range("K122").Value = 9356.625
$ws.Range("M122").Value = -6906.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2267.1316
$ws.Range("I132").Value = 1581.2693
$ws.Range("K132").Value = 4743.8079
$ws.Range("M132").Value = -2213.8079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 1875
$ws.Range("K122").Value = 5625
$ws.Range("M122").Value = -3175

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4027.2415
$ws.Range("J126").Value = 2958.9375
$ws.Range("L126").Value = 8876.8125
$ws.Range("N126").Value = -13816.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3596.652
$ws.Range("I132").Value = 3279.6216
$ws.Range("K132").Value = 9838.864799999999
$ws.Range("M132").Value = -7308.864799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1199.25
$ws.Range("I136").Value = 1081
$ws.Range("K136").Value = 3243
$ws.Range("M136").Value = -693
